# Add a 2022 data column (column S) to the "17.1.2" statistics table and
# move the active selection, matching the upstream gh-pages data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New year header for 2022, formatted like the existing R4 (2021) header cell.
$ws.Range("S4").Value = 2022
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)   # xlPasteFormats

# New data point for 2022, formatted like the existing R5 (2021) value cell.
$ws.Range("S5").Value = 76.1
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)   # xlPasteFormats

# Clear the marching-ants marquee left behind by Copy.
$excel.CutCopyMode = $false

# Update the active cell/selection as recorded in the saved view state.
$ws.Range("P8").Select() | Out-Null
